$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 792.46155
$ws.Range("I32").Value = 349.6
$ws.Range("J32").Value = 1069.25
$ws.Range("K32").Value = 349.6
$ws.Range("L32").Value = 1069.25
$ws.Range("M32").Value = -23.60000000000002
$ws.Range("N32").Value = -1721.25

$ws.Range("H75").Value = 19666.334
$ws.Range("J75").Value = 19666.334
$ws.Range("L75").Value = 19666.334
$ws.Range("N75").Value = -21538.334

$ws.Range("H78").Value = 19666.334
$ws.Range("J78").Value = 19666.334
$ws.Range("L78").Value = 58999.00199999999
$ws.Range("N78").Value = -68359.00199999999

$ws.Range("H98").Value = 4269
$ws.Range("I98").Value = 5448.0835
$ws.Range("J98").Value = 1439.2
$ws.Range("K98").Value = 5448.0835
$ws.Range("L98").Value = 1439.2
$ws.Range("M98").Value = -3950.0835
$ws.Range("N98").Value = -4435.2

$ws.Range("H122").Value = 4269
$ws.Range("I122").Value = 5448.0835
$ws.Range("J122").Value = 1439.2
$ws.Range("K122").Value = 16344.2505
$ws.Range("L122").Value = 4317.6
$ws.Range("M122").Value = -13894.2505
$ws.Range("N122").Value = -9217.6

$ws.Range("H125").Value = 3433.75
$ws.Range("I125").Value = 935
$ws.Range("K125").Value = 8415
$ws.Range("M125").Value = -5955

$ws.Range("H137").Value = 1424.629
$ws.Range("I137").Value = 1340.5714
$ws.Range("K137").Value = 4021.7142
$ws.Range("M137").Value = -1471.7142

$ws.Range("H138").Value = 485248.62
$ws.Range("I138").Value = 1359.3334
$ws.Range("J138").Value = 692629.75
$ws.Range("K138").Value = 4078.0002
$ws.Range("L138").Value = 2077889.25
$ws.Range("M138").Value = 1061.9998
$ws.Range("N138").Value = -2088169.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3689.46
$ws.Range("I32").Value = 3090.5557
$ws.Range("J32").Value = 9079.6
$ws.Range("K32").Value = 3090.5557
$ws.Range("L32").Value = 9079.6
$ws.Range("M32").Value = -2803.5557
$ws.Range("N32").Value = -9653.6

$ws.Range("H61").Value = 58824704
$ws.Range("I61").Value = 71429480
$ws.Range("K61").Value = 71429480
$ws.Range("M61").Value = -71429268

$ws.Range("H74").Value = 2358.3125
$ws.Range("I74").Value = 1986.6154
$ws.Range("J74").Value = 3969
$ws.Range("K74").Value = 1986.6154
$ws.Range("L74").Value = 3969
$ws.Range("M74").Value = -1112.6154
$ws.Range("N74").Value = -5717

$ws.Range("H77").Value = 2358.3125
$ws.Range("I77").Value = 1986.6154
$ws.Range("J77").Value = 3969
$ws.Range("K77").Value = 9933.076999999999
$ws.Range("L77").Value = 19845
$ws.Range("M77").Value = -5565.076999999999
$ws.Range("N77").Value = -28581

$ws.Range("H122").Value = 3523.8
$ws.Range("I122").Value = 3702.4
$ws.Range("J122").Value = 3345.2
$ws.Range("K122").Value = 11107.2
$ws.Range("L122").Value = 10035.6
$ws.Range("M122").Value = -8657.200000000001
$ws.Range("N122").Value = -14935.6

$ws.Range("H136").Value = 58824704
$ws.Range("I136").Value = 71429480
$ws.Range("K136").Value = 214288440
$ws.Range("M136").Value = -214285890

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 452
$ws.Range("I5").Value = 404
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 404
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -291
$ws.Range("N5").Value = -726

$ws.Range("H37").Value = 3833.3333
$ws.Range("I37").Value = 750
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 750
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -613
$ws.Range("N37").Value = -10274

$ws.Range("H134").Value = 1271.8422
$ws.Range("I134").Value = 1127.4117
$ws.Range("K134").Value = 3382.2351
$ws.Range("M134").Value = -847.2351000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1269.3191
$ws.Range("I31").Value = 1269.3191
$ws.Range("K31").Value = 1269.3191
$ws.Range("M31").Value = -974.3190999999999

$ws.Range("H34").Value = 1269.3191
$ws.Range("I34").Value = 1269.3191
$ws.Range("K34").Value = 1269.3191
$ws.Range("M34").Value = -1067.3191

$ws.Range("H51").Value = 23750
$ws.Range("I51").Value = 20000
$ws.Range("J51").Value = 25000
$ws.Range("K51").Value = 20000
$ws.Range("L51").Value = 25000
$ws.Range("M51").Value = -19264
$ws.Range("N51").Value = -26472

$ws.Range("H61").Value = 23750
$ws.Range("I61").Value = 20000
$ws.Range("J61").Value = 25000
$ws.Range("K61").Value = 20000
$ws.Range("L61").Value = 25000
$ws.Range("M61").Value = -19652
$ws.Range("N61").Value = -25696

$ws.Range("H94").Value = 2393.6365
$ws.Range("J94").Value = 2645.7144
$ws.Range("L94").Value = 2645.7144
$ws.Range("N94").Value = -3547.7144

$ws.Range("H111").Value = 50999.5
$ws.Range("J111").Value = 50999.5
$ws.Range("L111").Value = 50999.5
$ws.Range("N111").Value = -59179.5

$ws.Range("H141").Value = 481091.53
$ws.Range("J141").Value = 519849.16
$ws.Range("L141").Value = 519849.16
$ws.Range("N141").Value = -530209.1599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 18116832
$ws.Range("I129").Value = 47619496
$ws.Range("J129").Value = 5209417
$ws.Range("K129").Value = 142858488
$ws.Range("L129").Value = 15628251
$ws.Range("M129").Value = -142853488
$ws.Range("N129").Value = -15638251

$ws.Range("H131").Value = 19233532
$ws.Range("I131").Value = 90909544
$ws.Range("J131").Value = 3382.6584
$ws.Range("K131").Value = 272728632
$ws.Range("L131").Value = 10147.9752
$ws.Range("M131").Value = -272723592
$ws.Range("N131").Value = -20227.9752

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 20750
$ws.Range("J52").Value = 20750
$ws.Range("L52").Value = 20750
$ws.Range("N52").Value = -21268

$ws.Range("H132").Value = 2702.257
$ws.Range("I132").Value = 2936.2778
$ws.Range("J132").Value = 2454.4707
$ws.Range("K132").Value = 8808.8334
$ws.Range("L132").Value = 7363.4121
$ws.Range("M132").Value = -6278.8334
$ws.Range("N132").Value = -12423.4121

$ws.Range("H136").Value = 13652.952
$ws.Range("J136").Value = 13652.952
$ws.Range("L136").Value = 40958.856
$ws.Range("N136").Value = -46058.856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 5000
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5224

$ws.Range("H15").Value = 5000
$ws.Range("J15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("N15").Value = -5340

$ws.Range("H61").Value = 1053.6364
$ws.Range("I61").Value = 937.0625
$ws.Range("J61").Value = 1364.5
$ws.Range("K61").Value = 937.0625
$ws.Range("L61").Value = 1364.5
$ws.Range("M61").Value = -735.0625
$ws.Range("N61").Value = -1768.5

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H113").Value = 1053.6364
$ws.Range("I113").Value = 937.0625
$ws.Range("J113").Value = 1364.5
$ws.Range("K113").Value = 937.0625
$ws.Range("L113").Value = 1364.5
$ws.Range("M113").Value = 1232.9375
$ws.Range("N113").Value = -5704.5

$ws.Range("H133").Value = 46421.25
$ws.Range("J133").Value = 46421.25
$ws.Range("L133").Value = 46421.25
$ws.Range("N133").Value = -51481.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 14334967
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 17201760
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 17201760
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -17201984

$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H110").Value = 19000
$ws.Range("J110").Value = 19000
$ws.Range("L110").Value = 19000
$ws.Range("N110").Value = -27180

$ws.Range("H132").Value = 1862.4324
$ws.Range("I132").Value = 1648.8485
$ws.Range("K132").Value = 4946.5455
$ws.Range("M132").Value = -2416.5455
